# Junction_Flooding_286.xlsx edit:
#  - round the (existing) row-5 measurement values to 2 decimal places
#    ("custom accuracy")
#  - drop row 6 entirely (data trimmed to 1000 points upstream -> fewer
#    rows survive here)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update row 5 values to their 2-decimal "custom accuracy" values ---
$rowFive = @{
    "B5"  = 9
    "C5"  = 6.66
    "E5"  = 19.84
    "F5"  = 15.77
    "G5"  = 7.03
    "H5"  = 26.26
    "I5"  = 11.02
    "J5"  = 4.86
    "K5"  = 6.92
    "L5"  = 7.95
    "M5"  = 8.539999999999999
    "N5"  = 2.16
    "O5"  = 7.14
    "P5"  = 10
    "Q5"  = 6.22
    "R5"  = 0.68
    "S5"  = 0.38
    "T5"  = 101.37
    "U5"  = 19.85
    "V5"  = 6.59
    "W5"  = 13.11
    "Z5"  = 12.92
    "AA5" = 5.82
    "AB5" = 5.26
    "AC5" = 6.16
    "AD5" = 8.359999999999999
    "AE5" = 0.53
    "AF5" = 23.56
    "AH5" = 8.24
}

foreach ($addr in $rowFive.Keys) {
    $ws.Range($addr).Value = $rowFive[$addr]
}

# --- 2. Remove row 6 (and shrink the used range / dimension to A1:AH5) ---
$ws.Rows.Item(6).Delete()
